$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap paired-row data (rows that represent matches whose order flipped) ---
$ws.Cells.Item(2,2).Value = 6832493  # B2
$ws.Cells.Item(2,6).Value = "Al Wahda Abu Dhabi"  # F2
$ws.Cells.Item(2,7).Value = "Al Bataeh"  # G2
$ws.Cells.Item(2,8).Value = 1  # H2
$ws.Cells.Item(2,9).Value = 2  # I2
$ws.Cells.Item(2,11).Value = 1.25  # K2
$ws.Cells.Item(2,12).Value = 6  # L2
$ws.Cells.Item(2,13).Value = 8.5  # M2
$ws.Cells.Item(2,14).Value = 1.25  # N2
$ws.Cells.Item(2,15).Value = 6  # O2
$ws.Cells.Item(2,16).Value = 8.5  # P2
$ws.Cells.Item(2,17).Value = -1.75  # Q2
$ws.Cells.Item(2,18).Value = 1.875  # R2
$ws.Cells.Item(2,19).Value = 1.925  # S2
$ws.Cells.Item(2,20).Value = 3.25  # T2
$ws.Cells.Item(2,21).Value = 1.975  # U2
$ws.Cells.Item(2,22).Value = 1.825  # V2
$ws.Cells.Item(2,25).Value = 7.5  # Y2
$ws.Cells.Item(2,27).Value = 0.925  # AA2
$ws.Cells.Item(2,28).Value = -0.5  # AB2
$ws.Cells.Item(2,29).Value = 0.4125  # AC2
$ws.Cells.Item(3,2).Value = 6832491  # B3
$ws.Cells.Item(3,6).Value = "Al Ittihad Kalba"  # F3
$ws.Cells.Item(3,7).Value = "Sharjah SCC"  # G3
$ws.Cells.Item(3,8).Value = 3  # H3
$ws.Cells.Item(3,9).Value = 4  # I3
$ws.Cells.Item(3,11).Value = 3.8  # K3
$ws.Cells.Item(3,12).Value = 3.6  # L3
$ws.Cells.Item(3,13).Value = 1.8  # M3
$ws.Cells.Item(3,14).Value = 3.8  # N3
$ws.Cells.Item(3,15).Value = 3.6  # O3
$ws.Cells.Item(3,16).Value = 1.8  # P3
$ws.Cells.Item(3,17).Value = 0.5  # Q3
$ws.Cells.Item(3,18).Value = 2  # R3
$ws.Cells.Item(3,19).Value = 1.8  # S3
$ws.Cells.Item(3,20).Value = 2.75  # T3
$ws.Cells.Item(3,21).Value = 1.925  # U3
$ws.Cells.Item(3,22).Value = 1.875  # V3
$ws.Cells.Item(3,25).Value = 0.8  # Y3
$ws.Cells.Item(3,27).Value = 0.8  # AA3
$ws.Cells.Item(3,28).Value = 0.925  # AB3
$ws.Cells.Item(3,29).Value = -1  # AC3
$ws.Cells.Item(16,2).Value = 6832502  # B16
$ws.Cells.Item(16,6).Value = "Al Nasr SC"  # F16
$ws.Cells.Item(16,7).Value = "Sharjah SCC"  # G16
$ws.Cells.Item(16,9).Value = 1  # I16
$ws.Cells.Item(16,11).Value = 3.75  # K16
$ws.Cells.Item(16,12).Value = 3.6  # L16
$ws.Cells.Item(16,13).Value = 1.909  # M16
$ws.Cells.Item(16,15).Value = 3.8  # O16
$ws.Cells.Item(16,16).Value = 1.8  # P16
$ws.Cells.Item(16,17).Value = 0.5  # Q16
$ws.Cells.Item(16,18).Value = 2  # R16
$ws.Cells.Item(16,19).Value = 1.8  # S16
$ws.Cells.Item(16,20).Value = 3  # T16
$ws.Cells.Item(16,21).Value = 2  # U16
$ws.Cells.Item(16,22).Value = 1.8  # V16
$ws.Cells.Item(16,25).Value = 0.8  # Y16
$ws.Cells.Item(16,27).Value = 0.8  # AA16
$ws.Cells.Item(16,28).Value = -1  # AB16
$ws.Cells.Item(16,29).Value = 0.8  # AC16
$ws.Cells.Item(17,2).Value = 6832504  # B17
$ws.Cells.Item(17,6).Value = "Al Ittihad Kalba"  # F17
$ws.Cells.Item(17,7).Value = "Al Jazira SC"  # G17
$ws.Cells.Item(17,9).Value = 4  # I17
$ws.Cells.Item(17,11).Value = 3.6  # K17
$ws.Cells.Item(17,12).Value = 3.8  # L17
$ws.Cells.Item(17,13).Value = 1.85  # M17
$ws.Cells.Item(17,15).Value = 4  # O17
$ws.Cells.Item(17,16).Value = 1.727  # P17
$ws.Cells.Item(17,17).Value = 0.75  # Q17
$ws.Cells.Item(17,18).Value = 1.875  # R17
$ws.Cells.Item(17,19).Value = 1.925  # S17
$ws.Cells.Item(17,20).Value = 3.25  # T17
$ws.Cells.Item(17,21).Value = 1.95  # U17
$ws.Cells.Item(17,22).Value = 1.85  # V17
$ws.Cells.Item(17,25).Value = 0.7270000000000001  # Y17
$ws.Cells.Item(17,27).Value = 0.925  # AA17
$ws.Cells.Item(17,28).Value = 0.95  # AB17
$ws.Cells.Item(17,29).Value = -1  # AC17
$ws.Cells.Item(20,2).Value = 6832505  # B20
$ws.Cells.Item(20,6).Value = "Al Bataeh"  # F20
$ws.Cells.Item(20,7).Value = "Shabab Al Ahli Dubai"  # G20
$ws.Cells.Item(20,8).Value = 0  # H20
$ws.Cells.Item(20,9).Value = 2  # I20
$ws.Cells.Item(20,10).Value = "A"  # J20
$ws.Cells.Item(20,11).Value = 5.75  # K20
$ws.Cells.Item(20,13).Value = 1.4  # M20
$ws.Cells.Item(20,14).Value = 5.75  # N20
$ws.Cells.Item(20,15).Value = 5  # O20
$ws.Cells.Item(20,16).Value = 1.444  # P20
$ws.Cells.Item(20,17).Value = 1.25  # Q20
$ws.Cells.Item(20,18).Value = 1.875  # R20
$ws.Cells.Item(20,19).Value = 1.925  # S20
$ws.Cells.Item(20,20).Value = 3  # T20
$ws.Cells.Item(20,21).Value = 1.85  # U20
$ws.Cells.Item(20,22).Value = 1.95  # V20
$ws.Cells.Item(20,23).Value = -1  # W20
$ws.Cells.Item(20,25).Value = 0.444  # Y20
$ws.Cells.Item(20,27).Value = 0.925  # AA20
$ws.Cells.Item(20,29).Value = 0.95  # AC20
$ws.Cells.Item(21,2).Value = 6832506  # B21
$ws.Cells.Item(21,6).Value = "Al Wahda Abu Dhabi"  # F21
$ws.Cells.Item(21,7).Value = "Hatta Dubai"  # G21
$ws.Cells.Item(21,8).Value = 1  # H21
$ws.Cells.Item(21,9).Value = 0  # I21
$ws.Cells.Item(21,10).Value = "H"  # J21
$ws.Cells.Item(21,11).Value = 1.4  # K21
$ws.Cells.Item(21,13).Value = 6.5  # M21
$ws.Cells.Item(21,14).Value = 1.25  # N21
$ws.Cells.Item(21,15).Value = 5.5  # O21
$ws.Cells.Item(21,16).Value = 10  # P21
$ws.Cells.Item(21,17).Value = -1.75  # Q21
$ws.Cells.Item(21,18).Value = 1.9  # R21
$ws.Cells.Item(21,19).Value = 1.9  # S21
$ws.Cells.Item(21,20).Value = 3.25  # T21
$ws.Cells.Item(21,21).Value = 1.825  # U21
$ws.Cells.Item(21,22).Value = 1.975  # V21
$ws.Cells.Item(21,23).Value = 0.25  # W21
$ws.Cells.Item(21,25).Value = -1  # Y21
$ws.Cells.Item(21,27).Value = 0.8999999999999999  # AA21
$ws.Cells.Item(21,29).Value = 0.9750000000000001  # AC21
$ws.Cells.Item(30,2).Value = 6832514  # B30
$ws.Cells.Item(30,6).Value = "Al Ittihad Kalba"  # F30
$ws.Cells.Item(30,7).Value = "Hatta Dubai"  # G30
$ws.Cells.Item(30,8).Value = 2  # H30
$ws.Cells.Item(30,11).Value = 1.6  # K30
$ws.Cells.Item(30,12).Value = 4  # L30
$ws.Cells.Item(30,13).Value = 4.75  # M30
$ws.Cells.Item(30,14).Value = 1.444  # N30
$ws.Cells.Item(30,15).Value = 4.75  # O30
$ws.Cells.Item(30,16).Value = 5.75  # P30
$ws.Cells.Item(30,17).Value = -1.25  # Q30
$ws.Cells.Item(30,18).Value = 1.975  # R30
$ws.Cells.Item(30,19).Value = 1.825  # S30
$ws.Cells.Item(30,20).Value = 3.25  # T30
$ws.Cells.Item(30,21).Value = 1.975  # U30
$ws.Cells.Item(30,22).Value = 1.825  # V30
$ws.Cells.Item(30,23).Value = 0.444  # W30
$ws.Cells.Item(30,26).Value = -0.5  # Z30
$ws.Cells.Item(30,27).Value = 0.4125  # AA30
$ws.Cells.Item(30,28).Value = -0.5  # AB30
$ws.Cells.Item(30,29).Value = 0.4125  # AC30
$ws.Cells.Item(31,2).Value = 6832519  # B31
$ws.Cells.Item(31,6).Value = "Al Bataeh"  # F31
$ws.Cells.Item(31,7).Value = "Khor Fakkan"  # G31
$ws.Cells.Item(31,8).Value = 3  # H31
$ws.Cells.Item(31,11).Value = 2.375  # K31
$ws.Cells.Item(31,12).Value = 3.4  # L31
$ws.Cells.Item(31,13).Value = 2.6  # M31
$ws.Cells.Item(31,14).Value = 2.25  # N31
$ws.Cells.Item(31,15).Value = 3.4  # O31
$ws.Cells.Item(31,16).Value = 2.875  # P31
$ws.Cells.Item(31,17).Value = -0.25  # Q31
$ws.Cells.Item(31,18).Value = 2.025  # R31
$ws.Cells.Item(31,19).Value = 1.775  # S31
$ws.Cells.Item(31,20).Value = 2.75  # T31
$ws.Cells.Item(31,21).Value = 1.825  # U31
$ws.Cells.Item(31,22).Value = 1.975  # V31
$ws.Cells.Item(31,23).Value = 1.25  # W31
$ws.Cells.Item(31,26).Value = 1.025  # Z31
$ws.Cells.Item(31,27).Value = -1  # AA31
$ws.Cells.Item(31,28).Value = 0.825  # AB31
$ws.Cells.Item(31,29).Value = -1  # AC31
$ws.Cells.Item(39,2).Value = 6832522  # B39
$ws.Cells.Item(39,6).Value = "Shabab Al Ahli Dubai"  # F39
$ws.Cells.Item(39,7).Value = "Al Nasr SC"  # G39
$ws.Cells.Item(39,8).Value = 3  # H39
$ws.Cells.Item(39,9).Value = 3  # I39
$ws.Cells.Item(39,10).Value = "D"  # J39
$ws.Cells.Item(39,11).Value = 1.45  # K39
$ws.Cells.Item(39,12).Value = 4.75  # L39
$ws.Cells.Item(39,13).Value = 5.25  # M39
$ws.Cells.Item(39,14).Value = 1.4  # N39
$ws.Cells.Item(39,15).Value = 5  # O39
$ws.Cells.Item(39,16).Value = 5.5  # P39
$ws.Cells.Item(39,17).Value = -1.25  # Q39
$ws.Cells.Item(39,18).Value = 1.9  # R39
$ws.Cells.Item(39,19).Value = 1.9  # S39
$ws.Cells.Item(39,20).Value = 3.25  # T39
$ws.Cells.Item(39,21).Value = 1.925  # U39
$ws.Cells.Item(39,22).Value = 1.875  # V39
$ws.Cells.Item(39,23).Value = -1  # W39
$ws.Cells.Item(39,24).Value = 4  # X39
$ws.Cells.Item(39,26).Value = -1  # Z39
$ws.Cells.Item(39,27).Value = 0.8999999999999999  # AA39
$ws.Cells.Item(39,28).Value = 0.925  # AB39
$ws.Cells.Item(39,29).Value = -1  # AC39
$ws.Cells.Item(40,2).Value = 6832520  # B40
$ws.Cells.Item(40,6).Value = "Al Wahda Abu Dhabi"  # F40
$ws.Cells.Item(40,7).Value = "Al Ittihad Kalba"  # G40
$ws.Cells.Item(40,8).Value = 2  # H40
$ws.Cells.Item(40,9).Value = 1  # I40
$ws.Cells.Item(40,10).Value = "H"  # J40
$ws.Cells.Item(40,11).Value = 1.666  # K40
$ws.Cells.Item(40,12).Value = 4  # L40
$ws.Cells.Item(40,13).Value = 4.2  # M40
$ws.Cells.Item(40,14).Value = 1.666  # N40
$ws.Cells.Item(40,15).Value = 4  # O40
$ws.Cells.Item(40,16).Value = 4.2  # P40
$ws.Cells.Item(40,17).Value = -0.75  # Q40
$ws.Cells.Item(40,18).Value = 1.85  # R40
$ws.Cells.Item(40,19).Value = 1.95  # S40
$ws.Cells.Item(40,20).Value = 3  # T40
$ws.Cells.Item(40,21).Value = 1.775  # U40
$ws.Cells.Item(40,22).Value = 2.025  # V40
$ws.Cells.Item(40,23).Value = 0.6659999999999999  # W40
$ws.Cells.Item(40,24).Value = -1  # X40
$ws.Cells.Item(40,26).Value = 0.425  # Z40
$ws.Cells.Item(40,27).Value = -0.5  # AA40
$ws.Cells.Item(40,28).Value = 0  # AB40
$ws.Cells.Item(40,29).Value = 0  # AC40
$ws.Cells.Item(41,2).Value = 6832524  # B41
$ws.Cells.Item(41,6).Value = "Khor Fakkan"  # F41
$ws.Cells.Item(41,7).Value = "Al Jazira SC"  # G41
$ws.Cells.Item(41,8).Value = 4  # H41
$ws.Cells.Item(41,9).Value = 2  # I41
$ws.Cells.Item(41,11).Value = 5.25  # K41
$ws.Cells.Item(41,12).Value = 4.2  # L41
$ws.Cells.Item(41,13).Value = 1.533  # M41
$ws.Cells.Item(41,14).Value = 5.75  # N41
$ws.Cells.Item(41,15).Value = 4.5  # O41
$ws.Cells.Item(41,16).Value = 1.444  # P41
$ws.Cells.Item(41,17).Value = 1.25  # Q41
$ws.Cells.Item(41,18).Value = 1.95  # R41
$ws.Cells.Item(41,19).Value = 1.85  # S41
$ws.Cells.Item(41,21).Value = 2  # U41
$ws.Cells.Item(41,22).Value = 1.8  # V41
$ws.Cells.Item(41,23).Value = 4.75  # W41
$ws.Cells.Item(41,26).Value = 0.95  # Z41
$ws.Cells.Item(41,27).Value = -1  # AA41
$ws.Cells.Item(41,28).Value = 1  # AB41
$ws.Cells.Item(42,2).Value = 6832702  # B42
$ws.Cells.Item(42,6).Value = "Al Ain SCC"  # F42
$ws.Cells.Item(42,7).Value = "Emirates Club RAK"  # G42
$ws.Cells.Item(42,8).Value = 3  # H42
$ws.Cells.Item(42,9).Value = 1  # I42
$ws.Cells.Item(42,11).Value = 1.166  # K42
$ws.Cells.Item(42,12).Value = 7  # L42
$ws.Cells.Item(42,13).Value = 11  # M42
$ws.Cells.Item(42,14).Value = 1.125  # N42
$ws.Cells.Item(42,15).Value = 8  # O42
$ws.Cells.Item(42,16).Value = 13  # P42
$ws.Cells.Item(42,17).Value = -2.25  # Q42
$ws.Cells.Item(42,18).Value = 1.875  # R42
$ws.Cells.Item(42,19).Value = 1.925  # S42
$ws.Cells.Item(42,21).Value = 1.95  # U42
$ws.Cells.Item(42,22).Value = 1.85  # V42
$ws.Cells.Item(42,23).Value = 0.125  # W42
$ws.Cells.Item(42,26).Value = -0.5  # Z42
$ws.Cells.Item(42,27).Value = 0.4625  # AA42
$ws.Cells.Item(42,28).Value = 0.95  # AB42
$ws.Cells.Item(55,2).Value = 6832534  # B55
$ws.Cells.Item(55,6).Value = "Khor Fakkan"  # F55
$ws.Cells.Item(55,7).Value = "Ajman SCC"  # G55
$ws.Cells.Item(55,8).Value = 1  # H55
$ws.Cells.Item(55,9).Value = 2  # I55
$ws.Cells.Item(55,10).Value = "A"  # J55
$ws.Cells.Item(55,11).Value = 2.2  # K55
$ws.Cells.Item(55,12).Value = 3.5  # L55
$ws.Cells.Item(55,13).Value = 2.8  # M55
$ws.Cells.Item(55,14).Value = 2.1  # N55
$ws.Cells.Item(55,15).Value = 3.5  # O55
$ws.Cells.Item(55,16).Value = 3  # P55
$ws.Cells.Item(55,17).Value = -0.25  # Q55
$ws.Cells.Item(55,18).Value = 1.85  # R55
$ws.Cells.Item(55,19).Value = 1.95  # S55
$ws.Cells.Item(55,20).Value = 2.75  # T55
$ws.Cells.Item(55,21).Value = 1.85  # U55
$ws.Cells.Item(55,22).Value = 1.95  # V55
$ws.Cells.Item(55,24).Value = -1  # X55
$ws.Cells.Item(55,25).Value = 2  # Y55
$ws.Cells.Item(55,26).Value = -1  # Z55
$ws.Cells.Item(55,27).Value = 0.95  # AA55
$ws.Cells.Item(55,28).Value = 0.425  # AB55
$ws.Cells.Item(55,29).Value = -0.5  # AC55
$ws.Cells.Item(56,2).Value = 6832533  # B56
$ws.Cells.Item(56,6).Value = "Al Ittihad Kalba"  # F56
$ws.Cells.Item(56,7).Value = "Al Nasr SC"  # G56
$ws.Cells.Item(56,8).Value = 0  # H56
$ws.Cells.Item(56,9).Value = 0  # I56
$ws.Cells.Item(56,10).Value = "D"  # J56
$ws.Cells.Item(56,11).Value = 2.45  # K56
$ws.Cells.Item(56,12).Value = 3.6  # L56
$ws.Cells.Item(56,13).Value = 2.45  # M56
$ws.Cells.Item(56,14).Value = 2.625  # N56
$ws.Cells.Item(56,15).Value = 3.6  # O56
$ws.Cells.Item(56,16).Value = 2.25  # P56
$ws.Cells.Item(56,17).Value = 0.25  # Q56
$ws.Cells.Item(56,18).Value = 1.8  # R56
$ws.Cells.Item(56,19).Value = 2  # S56
$ws.Cells.Item(56,20).Value = 3  # T56
$ws.Cells.Item(56,21).Value = 1.925  # U56
$ws.Cells.Item(56,22).Value = 1.875  # V56
$ws.Cells.Item(56,24).Value = 2.6  # X56
$ws.Cells.Item(56,25).Value = -1  # Y56
$ws.Cells.Item(56,26).Value = 0.4  # Z56
$ws.Cells.Item(56,27).Value = -0.5  # AA56
$ws.Cells.Item(56,28).Value = -1  # AB56
$ws.Cells.Item(56,29).Value = 0.875  # AC56
$ws.Cells.Item(89,2).Value = 6832568  # B89
$ws.Cells.Item(89,6).Value = "Hatta Dubai"  # F89
$ws.Cells.Item(89,7).Value = "Khor Fakkan"  # G89
$ws.Cells.Item(89,11).Value = 2.9  # K89
$ws.Cells.Item(89,12).Value = 3.6  # L89
$ws.Cells.Item(89,13).Value = 2.15  # M89
$ws.Cells.Item(89,14).Value = 3.1  # N89
$ws.Cells.Item(89,15).Value = 3.5  # O89
$ws.Cells.Item(89,16).Value = 2.05  # P89
$ws.Cells.Item(89,17).Value = 0.25  # Q89
$ws.Cells.Item(89,18).Value = 1.95  # R89
$ws.Cells.Item(89,19).Value = 1.85  # S89
$ws.Cells.Item(89,20).Value = 2.75  # T89
$ws.Cells.Item(89,21).Value = 1.8  # U89
$ws.Cells.Item(89,22).Value = 2  # V89
$ws.Cells.Item(89,25).Value = 1.05  # Y89
$ws.Cells.Item(89,26).Value = -1  # Z89
$ws.Cells.Item(89,27).Value = 0.8500000000000001  # AA89
$ws.Cells.Item(89,29).Value = 1  # AC89
$ws.Cells.Item(90,2).Value = 6832569  # B90
$ws.Cells.Item(90,6).Value = "Emirates Club RAK"  # F90
$ws.Cells.Item(90,7).Value = "Al Wasl SC"  # G90
$ws.Cells.Item(90,11).Value = 7  # K90
$ws.Cells.Item(90,12).Value = 6.5  # L90
$ws.Cells.Item(90,13).Value = 1.285  # M90
$ws.Cells.Item(90,14).Value = 10  # N90
$ws.Cells.Item(90,15).Value = 9  # O90
$ws.Cells.Item(90,16).Value = 1.142  # P90
$ws.Cells.Item(90,17).Value = 2.5  # Q90
$ws.Cells.Item(90,18).Value = 1.8  # R90
$ws.Cells.Item(90,19).Value = 2  # S90
$ws.Cells.Item(90,20).Value = 4  # T90
$ws.Cells.Item(90,21).Value = 1.95  # U90
$ws.Cells.Item(90,22).Value = 1.85  # V90
$ws.Cells.Item(90,25).Value = 0.1419999999999999  # Y90
$ws.Cells.Item(90,26).Value = 0.8  # Z90
$ws.Cells.Item(90,27).Value = -1  # AA90
$ws.Cells.Item(90,29).Value = 0.8500000000000001  # AC90
$ws.Cells.Item(98,2).Value = 6832572  # B98
$ws.Cells.Item(98,6).Value = "Khor Fakkan"  # F98
$ws.Cells.Item(98,7).Value = "Shabab Al Ahli Dubai"  # G98
$ws.Cells.Item(98,8).Value = 0  # H98
$ws.Cells.Item(98,9).Value = 3  # I98
$ws.Cells.Item(98,10).Value = "A"  # J98
$ws.Cells.Item(98,11).Value = 6  # K98
$ws.Cells.Item(98,12).Value = 5.5  # L98
$ws.Cells.Item(98,13).Value = 1.333  # M98
$ws.Cells.Item(98,14).Value = 6.5  # N98
$ws.Cells.Item(98,15).Value = 5.75  # O98
$ws.Cells.Item(98,16).Value = 1.3  # P98
$ws.Cells.Item(98,17).Value = 1.75  # Q98
$ws.Cells.Item(98,18).Value = 1.8  # R98
$ws.Cells.Item(98,19).Value = 2  # S98
$ws.Cells.Item(98,20).Value = 3.5  # T98
$ws.Cells.Item(98,21).Value = 1.8  # U98
$ws.Cells.Item(98,22).Value = 2  # V98
$ws.Cells.Item(98,23).Value = -1  # W98
$ws.Cells.Item(98,25).Value = 0.3  # Y98
$ws.Cells.Item(98,26).Value = -1  # Z98
$ws.Cells.Item(98,27).Value = 1  # AA98
$ws.Cells.Item(98,29).Value = 1  # AC98
$ws.Cells.Item(99,2).Value = 6832576  # B99
$ws.Cells.Item(99,6).Value = "Ajman SCC"  # F99
$ws.Cells.Item(99,7).Value = "Emirates Club RAK"  # G99
$ws.Cells.Item(99,8).Value = 2  # H99
$ws.Cells.Item(99,9).Value = 0  # I99
$ws.Cells.Item(99,10).Value = "H"  # J99
$ws.Cells.Item(99,11).Value = 1.6  # K99
$ws.Cells.Item(99,12).Value = 4  # L99
$ws.Cells.Item(99,13).Value = 4.75  # M99
$ws.Cells.Item(99,14).Value = 1.7  # N99
$ws.Cells.Item(99,15).Value = 3.8  # O99
$ws.Cells.Item(99,16).Value = 4.333  # P99
$ws.Cells.Item(99,17).Value = -0.75  # Q99
$ws.Cells.Item(99,18).Value = 1.875  # R99
$ws.Cells.Item(99,19).Value = 1.925  # S99
$ws.Cells.Item(99,20).Value = 3.25  # T99
$ws.Cells.Item(99,21).Value = 1.975  # U99
$ws.Cells.Item(99,22).Value = 1.825  # V99
$ws.Cells.Item(99,23).Value = 0.7  # W99
$ws.Cells.Item(99,25).Value = -1  # Y99
$ws.Cells.Item(99,26).Value = 0.875  # Z99
$ws.Cells.Item(99,27).Value = -1  # AA99
$ws.Cells.Item(99,29).Value = 0.825  # AC99
$ws.Cells.Item(103,2).Value = 6832580  # B103
$ws.Cells.Item(103,6).Value = "Hatta Dubai"  # F103
$ws.Cells.Item(103,7).Value = "Al Wahda Abu Dhabi"  # G103
$ws.Cells.Item(103,9).Value = 4  # I103
$ws.Cells.Item(103,11).Value = 9  # K103
$ws.Cells.Item(103,12).Value = 5.5  # L103
$ws.Cells.Item(103,13).Value = 1.25  # M103
$ws.Cells.Item(103,14).Value = 6.5  # N103
$ws.Cells.Item(103,16).Value = 1.4  # P103
$ws.Cells.Item(103,17).Value = 1.25  # Q103
$ws.Cells.Item(103,20).Value = 3.25  # T103
$ws.Cells.Item(103,21).Value = 1.975  # U103
$ws.Cells.Item(103,22).Value = 1.825  # V103
$ws.Cells.Item(103,25).Value = 0.3999999999999999  # Y103
$ws.Cells.Item(103,28).Value = 0.9750000000000001  # AB103
$ws.Cells.Item(103,29).Value = -1  # AC103
$ws.Cells.Item(104,2).Value = 6832583  # B104
$ws.Cells.Item(104,6).Value = "Baniyas SC"  # F104
$ws.Cells.Item(104,7).Value = "Al Wasl SC"  # G104
$ws.Cells.Item(104,9).Value = 2  # I104
$ws.Cells.Item(104,11).Value = 5.25  # K104
$ws.Cells.Item(104,12).Value = 5.25  # L104
$ws.Cells.Item(104,13).Value = 1.4  # M104
$ws.Cells.Item(104,14).Value = 4.75  # N104
$ws.Cells.Item(104,16).Value = 1.5  # P104
$ws.Cells.Item(104,17).Value = 1  # Q104
$ws.Cells.Item(104,20).Value = 3  # T104
$ws.Cells.Item(104,21).Value = 1.825  # U104
$ws.Cells.Item(104,22).Value = 1.975  # V104
$ws.Cells.Item(104,25).Value = 0.5  # Y104
$ws.Cells.Item(104,28).Value = -1  # AB104
$ws.Cells.Item(104,29).Value = 0.9750000000000001  # AC104
$ws.Cells.Item(106,2).Value = 6832577  # B106
$ws.Cells.Item(106,6).Value = "Shabab Al Ahli Dubai"  # F106
$ws.Cells.Item(106,7).Value = "Al Bataeh"  # G106
$ws.Cells.Item(106,9).Value = 1  # I106
$ws.Cells.Item(106,10).Value = "H"  # J106
$ws.Cells.Item(106,11).Value = 1.333  # K106
$ws.Cells.Item(106,12).Value = 5  # L106
$ws.Cells.Item(106,13).Value = 7.5  # M106
$ws.Cells.Item(106,14).Value = 1.25  # N106
$ws.Cells.Item(106,15).Value = 5.5  # O106
$ws.Cells.Item(106,16).Value = 8.5  # P106
$ws.Cells.Item(106,17).Value = -1.75  # Q106
$ws.Cells.Item(106,18).Value = 1.9  # R106
$ws.Cells.Item(106,19).Value = 1.9  # S106
$ws.Cells.Item(106,21).Value = 1.875  # U106
$ws.Cells.Item(106,22).Value = 1.925  # V106
$ws.Cells.Item(106,23).Value = 0.25  # W106
$ws.Cells.Item(106,24).Value = -1  # X106
$ws.Cells.Item(106,27).Value = 0.8999999999999999  # AA106
$ws.Cells.Item(106,28).Value = -1  # AB106
$ws.Cells.Item(106,29).Value = 0.925  # AC106
$ws.Cells.Item(107,2).Value = 6832582  # B107
$ws.Cells.Item(107,6).Value = "Al Jazira SC"  # F107
$ws.Cells.Item(107,7).Value = "Al Ittihad Kalba"  # G107
$ws.Cells.Item(107,9).Value = 2  # I107
$ws.Cells.Item(107,10).Value = "D"  # J107
$ws.Cells.Item(107,11).Value = 1.7  # K107
$ws.Cells.Item(107,12).Value = 4.2  # L107
$ws.Cells.Item(107,13).Value = 4  # M107
$ws.Cells.Item(107,14).Value = 1.8  # N107
$ws.Cells.Item(107,15).Value = 4  # O107
$ws.Cells.Item(107,16).Value = 3.5  # P107
$ws.Cells.Item(107,17).Value = -0.5  # Q107
$ws.Cells.Item(107,18).Value = 1.8  # R107
$ws.Cells.Item(107,19).Value = 2  # S107
$ws.Cells.Item(107,21).Value = 1.9  # U107
$ws.Cells.Item(107,22).Value = 1.9  # V107
$ws.Cells.Item(107,23).Value = -1  # W107
$ws.Cells.Item(107,24).Value = 3  # X107
$ws.Cells.Item(107,27).Value = 1  # AA107
$ws.Cells.Item(107,28).Value = 0.8999999999999999  # AB107
$ws.Cells.Item(107,29).Value = -1  # AC107
$ws.Cells.Item(110,2).Value = 6832714  # B110
$ws.Cells.Item(110,6).Value = "Emirates Club RAK"  # F110
$ws.Cells.Item(110,7).Value = "Baniyas SC"  # G110
$ws.Cells.Item(110,9).Value = 2  # I110
$ws.Cells.Item(110,10).Value = "A"  # J110
$ws.Cells.Item(110,11).Value = 3.4  # K110
$ws.Cells.Item(110,12).Value = 3.6  # L110
$ws.Cells.Item(110,13).Value = 1.909  # M110
$ws.Cells.Item(110,14).Value = 3.3  # N110
$ws.Cells.Item(110,15).Value = 3.75  # O110
$ws.Cells.Item(110,16).Value = 1.909  # P110
$ws.Cells.Item(110,17).Value = 0.5  # Q110
$ws.Cells.Item(110,21).Value = 1.85  # U110
$ws.Cells.Item(110,22).Value = 1.95  # V110
$ws.Cells.Item(110,23).Value = -1  # W110
$ws.Cells.Item(110,25).Value = 0.909  # Y110
$ws.Cells.Item(110,28).Value = -0.5  # AB110
$ws.Cells.Item(110,29).Value = 0.475  # AC110
$ws.Cells.Item(111,2).Value = 6832584  # B111
$ws.Cells.Item(111,6).Value = "Al Nasr SC"  # F111
$ws.Cells.Item(111,7).Value = "Hatta Dubai"  # G111
$ws.Cells.Item(111,9).Value = 0  # I111
$ws.Cells.Item(111,10).Value = "H"  # J111
$ws.Cells.Item(111,11).Value = 1.363  # K111
$ws.Cells.Item(111,12).Value = 5  # L111
$ws.Cells.Item(111,13).Value = 6  # M111
$ws.Cells.Item(111,14).Value = 1.25  # N111
$ws.Cells.Item(111,15).Value = 6  # O111
$ws.Cells.Item(111,16).Value = 8  # P111
$ws.Cells.Item(111,17).Value = -1.75  # Q111
$ws.Cells.Item(111,21).Value = 1.825  # U111
$ws.Cells.Item(111,22).Value = 1.975  # V111
$ws.Cells.Item(111,23).Value = 0.25  # W111
$ws.Cells.Item(111,25).Value = -1  # Y111
$ws.Cells.Item(111,28).Value = -1  # AB111
$ws.Cells.Item(111,29).Value = 0.9750000000000001  # AC111
$ws.Cells.Item(112,2).Value = 6832587  # B112
$ws.Cells.Item(112,6).Value = "Al Ittihad Kalba"  # F112
$ws.Cells.Item(112,7).Value = "Ajman SCC"  # G112
$ws.Cells.Item(112,8).Value = 2  # H112
$ws.Cells.Item(112,9).Value = 4  # I112
$ws.Cells.Item(112,10).Value = "A"  # J112
$ws.Cells.Item(112,11).Value = 1.85  # K112
$ws.Cells.Item(112,12).Value = 3.75  # L112
$ws.Cells.Item(112,13).Value = 3.5  # M112
$ws.Cells.Item(112,14).Value = 1.95  # N112
$ws.Cells.Item(112,15).Value = 3.6  # O112
$ws.Cells.Item(112,16).Value = 3.25  # P112
$ws.Cells.Item(112,17).Value = -0.25  # Q112
$ws.Cells.Item(112,18).Value = 1.8  # R112
$ws.Cells.Item(112,19).Value = 2  # S112
$ws.Cells.Item(112,20).Value = 3  # T112
$ws.Cells.Item(112,21).Value = 1.85  # U112
$ws.Cells.Item(112,22).Value = 1.95  # V112
$ws.Cells.Item(112,24).Value = -1  # X112
$ws.Cells.Item(112,25).Value = 2.25  # Y112
$ws.Cells.Item(112,26).Value = -1  # Z112
$ws.Cells.Item(112,27).Value = 1  # AA112
$ws.Cells.Item(112,28).Value = 0.8500000000000001  # AB112
$ws.Cells.Item(113,2).Value = 6832586  # B113
$ws.Cells.Item(113,6).Value = "Al Wahda Abu Dhabi"  # F113
$ws.Cells.Item(113,7).Value = "Shabab Al Ahli Dubai"  # G113
$ws.Cells.Item(113,8).Value = 3  # H113
$ws.Cells.Item(113,9).Value = 3  # I113
$ws.Cells.Item(113,10).Value = "D"  # J113
$ws.Cells.Item(113,11).Value = 2.6  # K113
$ws.Cells.Item(113,12).Value = 3.5  # L113
$ws.Cells.Item(113,13).Value = 2.375  # M113
$ws.Cells.Item(113,14).Value = 2.55  # N113
$ws.Cells.Item(113,15).Value = 3.3  # O113
$ws.Cells.Item(113,16).Value = 2.45  # P113
$ws.Cells.Item(113,17).Value = 0  # Q113
$ws.Cells.Item(113,18).Value = 1.975  # R113
$ws.Cells.Item(113,19).Value = 1.825  # S113
$ws.Cells.Item(113,20).Value = 2.75  # T113
$ws.Cells.Item(113,21).Value = 1.9  # U113
$ws.Cells.Item(113,22).Value = 1.9  # V113
$ws.Cells.Item(113,24).Value = 2.3  # X113
$ws.Cells.Item(113,25).Value = -1  # Y113
$ws.Cells.Item(113,26).Value = 0  # Z113
$ws.Cells.Item(113,27).Value = 0  # AA113
$ws.Cells.Item(113,28).Value = 0.8999999999999999  # AB113

# --- Append new rows 117-123 ---

# Row 117
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(117,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(117,5))
$ws.Cells.Item(117,1).Value = 115  # A117
$ws.Cells.Item(117,2).Value = 6832594  # B117
$ws.Cells.Item(117,6).Value = "Hatta Dubai"  # F117
$ws.Cells.Item(117,7).Value = "Al Ittihad Kalba"  # G117
$ws.Cells.Item(117,11).Value = 5  # K117
$ws.Cells.Item(117,12).Value = 4  # L117
$ws.Cells.Item(117,13).Value = 1.533  # M117
$ws.Cells.Item(117,14).Value = 5  # N117
$ws.Cells.Item(117,15).Value = 4  # O117
$ws.Cells.Item(117,16).Value = 1.533  # P117
$ws.Cells.Item(117,17).Value = 1  # Q117
$ws.Cells.Item(117,18).Value = 1.925  # R117
$ws.Cells.Item(117,19).Value = 1.875  # S117
$ws.Cells.Item(117,20).Value = 3.25  # T117
$ws.Cells.Item(117,21).Value = 2  # U117
$ws.Cells.Item(117,22).Value = 1.8  # V117
$ws.Cells.Item(117,23).Value = 0  # W117
$ws.Cells.Item(117,24).Value = 0  # X117
$ws.Cells.Item(117,25).Value = 0  # Y117
$ws.Cells.Item(117,26).Value = 0  # Z117
$ws.Cells.Item(117,27).Value = 0  # AA117
$ws.Cells.Item(117,3).Value = "UAE Premier League"  # C117
$ws.Cells.Item(117,4).Value = "UAE Premier League"  # D117
$ws.Cells.Item(117,5).Value = 45388.625  # E117

# Row 118
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(118,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(118,5))
$ws.Cells.Item(118,1).Value = 116  # A118
$ws.Cells.Item(118,2).Value = 6832592  # B118
$ws.Cells.Item(118,6).Value = "Al Wahda Abu Dhabi"  # F118
$ws.Cells.Item(118,7).Value = "Baniyas SC"  # G118
$ws.Cells.Item(118,11).Value = 1.363  # K118
$ws.Cells.Item(118,12).Value = 4.75  # L118
$ws.Cells.Item(118,13).Value = 6.5  # M118
$ws.Cells.Item(118,14).Value = 1.363  # N118
$ws.Cells.Item(118,15).Value = 4.75  # O118
$ws.Cells.Item(118,16).Value = 6.5  # P118
$ws.Cells.Item(118,17).Value = -1.5  # Q118
$ws.Cells.Item(118,18).Value = 1.975  # R118
$ws.Cells.Item(118,19).Value = 1.825  # S118
$ws.Cells.Item(118,20).Value = 3.25  # T118
$ws.Cells.Item(118,21).Value = 1.9  # U118
$ws.Cells.Item(118,22).Value = 1.9  # V118
$ws.Cells.Item(118,23).Value = 0  # W118
$ws.Cells.Item(118,24).Value = 0  # X118
$ws.Cells.Item(118,25).Value = 0  # Y118
$ws.Cells.Item(118,26).Value = 0  # Z118
$ws.Cells.Item(118,27).Value = 0  # AA118
$ws.Cells.Item(118,3).Value = "UAE Premier League"  # C118
$ws.Cells.Item(118,4).Value = "UAE Premier League"  # D118
$ws.Cells.Item(118,5).Value = 45388.625  # E118

# Row 119
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(119,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(119,5))
$ws.Cells.Item(119,1).Value = 117  # A119
$ws.Cells.Item(119,2).Value = 6832593  # B119
$ws.Cells.Item(119,6).Value = "Khor Fakkan"  # F119
$ws.Cells.Item(119,7).Value = "Al Bataeh"  # G119
$ws.Cells.Item(119,11).Value = 2.6  # K119
$ws.Cells.Item(119,12).Value = 3.4  # L119
$ws.Cells.Item(119,13).Value = 2.45  # M119
$ws.Cells.Item(119,14).Value = 2.6  # N119
$ws.Cells.Item(119,15).Value = 3.4  # O119
$ws.Cells.Item(119,16).Value = 2.45  # P119
$ws.Cells.Item(119,17).Value = 0  # Q119
$ws.Cells.Item(119,18).Value = 1.925  # R119
$ws.Cells.Item(119,19).Value = 1.875  # S119
$ws.Cells.Item(119,20).Value = 3  # T119
$ws.Cells.Item(119,21).Value = 1.85  # U119
$ws.Cells.Item(119,22).Value = 1.95  # V119
$ws.Cells.Item(119,23).Value = 0  # W119
$ws.Cells.Item(119,24).Value = 0  # X119
$ws.Cells.Item(119,25).Value = 0  # Y119
$ws.Cells.Item(119,26).Value = 0  # Z119
$ws.Cells.Item(119,27).Value = 0  # AA119
$ws.Cells.Item(119,3).Value = "UAE Premier League"  # C119
$ws.Cells.Item(119,4).Value = "UAE Premier League"  # D119
$ws.Cells.Item(119,5).Value = 45388.625  # E119

# Row 120
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(120,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(120,5))
$ws.Cells.Item(120,1).Value = 118  # A120
$ws.Cells.Item(120,2).Value = 6832589  # B120
$ws.Cells.Item(120,6).Value = "Shabab Al Ahli Dubai"  # F120
$ws.Cells.Item(120,7).Value = "Al Jazira SC"  # G120
$ws.Cells.Item(120,11).Value = 1.666  # K120
$ws.Cells.Item(120,12).Value = 4  # L120
$ws.Cells.Item(120,13).Value = 4  # M120
$ws.Cells.Item(120,14).Value = 1.666  # N120
$ws.Cells.Item(120,15).Value = 4  # O120
$ws.Cells.Item(120,16).Value = 4  # P120
$ws.Cells.Item(120,17).Value = -0.75  # Q120
$ws.Cells.Item(120,18).Value = 1.875  # R120
$ws.Cells.Item(120,19).Value = 1.925  # S120
$ws.Cells.Item(120,20).Value = 3.5  # T120
$ws.Cells.Item(120,21).Value = 1.925  # U120
$ws.Cells.Item(120,22).Value = 1.875  # V120
$ws.Cells.Item(120,23).Value = 0  # W120
$ws.Cells.Item(120,24).Value = 0  # X120
$ws.Cells.Item(120,25).Value = 0  # Y120
$ws.Cells.Item(120,26).Value = 0  # Z120
$ws.Cells.Item(120,27).Value = 0  # AA120
$ws.Cells.Item(120,3).Value = "UAE Premier League"  # C120
$ws.Cells.Item(120,4).Value = "UAE Premier League"  # D120
$ws.Cells.Item(120,5).Value = 45389.625  # E120

# Row 121
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(121,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(121,5))
$ws.Cells.Item(121,1).Value = 119  # A121
$ws.Cells.Item(121,2).Value = 6832591  # B121
$ws.Cells.Item(121,6).Value = "Al Ain SCC"  # F121
$ws.Cells.Item(121,7).Value = "Al Nasr SC"  # G121
$ws.Cells.Item(121,11).Value = 1.571  # K121
$ws.Cells.Item(121,12).Value = 4.333  # L121
$ws.Cells.Item(121,13).Value = 4.5  # M121
$ws.Cells.Item(121,14).Value = 1.571  # N121
$ws.Cells.Item(121,15).Value = 4.333  # O121
$ws.Cells.Item(121,16).Value = 4.5  # P121
$ws.Cells.Item(121,17).Value = -1  # Q121
$ws.Cells.Item(121,18).Value = 2  # R121
$ws.Cells.Item(121,19).Value = 1.8  # S121
$ws.Cells.Item(121,20).Value = 3  # T121
$ws.Cells.Item(121,21).Value = 1.8  # U121
$ws.Cells.Item(121,22).Value = 2  # V121
$ws.Cells.Item(121,23).Value = 0  # W121
$ws.Cells.Item(121,24).Value = 0  # X121
$ws.Cells.Item(121,25).Value = 0  # Y121
$ws.Cells.Item(121,26).Value = 0  # Z121
$ws.Cells.Item(121,27).Value = 0  # AA121
$ws.Cells.Item(121,3).Value = "UAE Premier League"  # C121
$ws.Cells.Item(121,4).Value = "UAE Premier League"  # D121
$ws.Cells.Item(121,5).Value = 45389.625  # E121

# Row 122
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(122,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(122,5))
$ws.Cells.Item(122,1).Value = 120  # A122
$ws.Cells.Item(122,2).Value = 8045291  # B122
$ws.Cells.Item(122,6).Value = "Ajman SCC"  # F122
$ws.Cells.Item(122,7).Value = "Al Wasl SC"  # G122
$ws.Cells.Item(122,11).Value = 6.5  # K122
$ws.Cells.Item(122,12).Value = 4.5  # L122
$ws.Cells.Item(122,13).Value = 1.4  # M122
$ws.Cells.Item(122,14).Value = 7.5  # N122
$ws.Cells.Item(122,15).Value = 5  # O122
$ws.Cells.Item(122,16).Value = 1.333  # P122
$ws.Cells.Item(122,17).Value = 1.5  # Q122
$ws.Cells.Item(122,18).Value = 1.875  # R122
$ws.Cells.Item(122,19).Value = 1.925  # S122
$ws.Cells.Item(122,20).Value = 3.25  # T122
$ws.Cells.Item(122,21).Value = 1.95  # U122
$ws.Cells.Item(122,22).Value = 1.85  # V122
$ws.Cells.Item(122,23).Value = 0  # W122
$ws.Cells.Item(122,24).Value = 0  # X122
$ws.Cells.Item(122,25).Value = 0  # Y122
$ws.Cells.Item(122,26).Value = 0  # Z122
$ws.Cells.Item(122,27).Value = 0  # AA122
$ws.Cells.Item(122,3).Value = "UAE Premier League"  # C122
$ws.Cells.Item(122,4).Value = "UAE Premier League"  # D122
$ws.Cells.Item(122,5).Value = 45390.625  # E122

# Row 123
$ws.Cells.Item(2,1).Copy($ws.Cells.Item(123,1))
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(123,5))
$ws.Cells.Item(123,1).Value = 121  # A123
$ws.Cells.Item(123,2).Value = 6832590  # B123
$ws.Cells.Item(123,6).Value = "Sharjah SCC"  # F123
$ws.Cells.Item(123,7).Value = "Emirates Club RAK"  # G123
$ws.Cells.Item(123,11).Value = 1.25  # K123
$ws.Cells.Item(123,12).Value = 5.5  # L123
$ws.Cells.Item(123,13).Value = 9.5  # M123
$ws.Cells.Item(123,14).Value = 1.25  # N123
$ws.Cells.Item(123,15).Value = 5.5  # O123
$ws.Cells.Item(123,16).Value = 9.5  # P123
$ws.Cells.Item(123,17).Value = -1.75  # Q123
$ws.Cells.Item(123,18).Value = 1.825  # R123
$ws.Cells.Item(123,19).Value = 1.975  # S123
$ws.Cells.Item(123,20).Value = 3.5  # T123
$ws.Cells.Item(123,21).Value = 2  # U123
$ws.Cells.Item(123,22).Value = 1.8  # V123
$ws.Cells.Item(123,23).Value = 0  # W123
$ws.Cells.Item(123,24).Value = 0  # X123
$ws.Cells.Item(123,25).Value = 0  # Y123
$ws.Cells.Item(123,26).Value = 0  # Z123
$ws.Cells.Item(123,27).Value = 0  # AA123
$ws.Cells.Item(123,3).Value = "UAE Premier League"  # C123
$ws.Cells.Item(123,4).Value = "UAE Premier League"  # D123
$ws.Cells.Item(123,5).Value = 45390.625  # E123